$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.851.21'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.642.51'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.66'
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.250'
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0622'
$ws.Range("E9").Value = '  -1.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.24'
$ws.Range("E10").Value = '  +0.49%  '
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.871.53'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.647.51'
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("E14").Value = '  -0.05%  '
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.39'
$ws.Range("E16").Value = '  +1.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.861.28'
$ws.Range("E17").Value = '  +0.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0729'
$ws.Range("E18").Value = '  -1.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.87'
$ws.Range("E19").Value = '  +0.89%  '
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.64'
$ws.Range("E21").Value = '  +6.08%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.37'
$ws.Range("E23").Value = '  -1.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.19'
$ws.Range("E24").Value = '  -1.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.66'
$ws.Range("E25").Value = '  +1.91%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.22'
$ws.Range("E28").Value = '  +1.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.71'
$ws.Range("E29").Value = '  +0.52%  '
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("E31").Value = '  +1.50%  '
$ws.Range("E32").Value = '  +1.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.99'
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.283.14'
$ws.Range("E34").Value = '  -0.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("E35").Value = '  +0.51%  '
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("E37").Value = '  -0.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.533'
$ws.Range("E38").Value = '  -0.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.818'
$ws.Range("E39").Value = '  -0.92%  '
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.34'
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.782.65'
$ws.Range("E43").Value = '  -0.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.11'
$ws.Range("E44").Value = '  -5.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.67'
$ws.Range("E45").Value = '  +1.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.96'
$ws.Range("E46").Value = '  -0.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.59'
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.58'
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0966'
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("E51").Value = '  -0.09%  '
